# Update the "counter" column (D) for bob and charlie rows per the commit's
# user-data spreadsheet update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 15
$ws.Range("D4").Value = 6
